# Updated cryptos list on Fri Oct 27 10:42:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.102.63'
$ws.Cells.Item(2, 5).Value = '  -0.55%  '
$ws.Cells.Item(3, 4).Value = '1.784.21'
$ws.Cells.Item(3, 5).Value = '  -2.64%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '225.00'
$ws.Cells.Item(5, 5).Value = '  -0.09%  '
$ws.Cells.Item(6, 5).Value = '  -1.50%  '
$ws.Cells.Item(7, 5).Value = '  +0.03%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '32.69'
$ws.Cells.Item(8, 5).Value = '  +2.17%  '
$ws.Cells.Item(9, 5).Value = '  -2.23%  '
$ws.Cells.Item(10, 5).Value = '  -1.56%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0937'
$ws.Cells.Item(11, 5).Value = '  +0.57%  '
$ws.Cells.Item(12, 4).Value = '2.041.12'
$ws.Cells.Item(12, 5).Value = '  -2.86%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.798.50'
$ws.Cells.Item(13, 5).Value = '  -1.92%  '
$ws.Cells.Item(14, 2).Value = 'Chainlink'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '10.98'
$ws.Cells.Item(14, 5).Value = '  +1.54%  '
$ws.Cells.Item(15, 4).Value = '34.029.63'
$ws.Cells.Item(15, 5).Value = '  -0.96%  '
$ws.Cells.Item(16, 5).Value = '  -3.71%  '
$ws.Cells.Item(17, 5).Value = '  -4.49%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '67.72'
$ws.Cells.Item(18, 5).Value = '  -3.02%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '245.11'
$ws.Cells.Item(19, 5).Value = '  -2.88%  '
$ws.Cells.Item(20, 5).Value = '  -0.73%  '
$ws.Cells.Item(21, 5).Value = '  +0.11%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.79'
$ws.Cells.Item(22, 5).Value = '  -3.62%  '
$ws.Cells.Item(23, 5).Value = '  -4.49%  '
$ws.Cells.Item(24, 5).Value = '  -2.66%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '160.20'
$ws.Cells.Item(25, 5).Value = '  -0.14%  '
$ws.Cells.Item(26, 5).Value = '  -2.20%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.04'
$ws.Cells.Item(27, 5).Value = '  -3.07%  '
$ws.Cells.Item(28, 5).Value = '  -2.47%  '
$ws.Cells.Item(29, 5).Value = '  +0.05%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.22'
$ws.Cells.Item(30, 5).Value = '  +0.11%  '
$ws.Cells.Item(31, 5).Value = '  -4.79%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.63'
$ws.Cells.Item(32, 5).Value = '  -4.35%  '
$ws.Cells.Item(33, 5).Value = '  -2.13%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.81'
$ws.Cells.Item(34, 5).Value = '  -5.26%  '
$ws.Cells.Item(35, 4).Value = '1.391.87'
$ws.Cells.Item(35, 5).Value = '  -3.82%  '
$ws.Cells.Item(36, 5).Value = '  -0.41%  '
$ws.Cells.Item(37, 5).Value = '  -1.99%  '
$ws.Cells.Item(38, 5).Value = '  -2.68%  '
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.20'
$ws.Cells.Item(39, 5).Value = '  +2.76%  '
$ws.Cells.Item(40, 2).Value = 'HuobiToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.35'
$ws.Cells.Item(40, 5).Value = '  -0.21%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.914'
$ws.Cells.Item(41, 5).Value = '  -5.30%  '
$ws.Cells.Item(42, 5).Value = '  -4.95%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '78.08'
$ws.Cells.Item(43, 5).Value = '  -4.70%  '
$ws.Cells.Item(44, 4).Value = '0.0₆0141'
$ws.Cells.Item(44, 5).Value = '  +12.54%  '
$ws.Cells.Item(45, 5).Value = '  +2.79%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '108.20'
$ws.Cells.Item(46, 5).Value = '  +1.43%  '
$ws.Cells.Item(47, 5).Value = '  -0.51%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '12.46'
$ws.Cells.Item(48, 5).Value = '  +4.75%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '5.83'
$ws.Cells.Item(49, 5).Value = '  -4.48%  '
$ws.Cells.Item(50, 4).Value = '1.940.62'
$ws.Cells.Item(50, 5).Value = '  -2.77%  '
$ws.Cells.Item(51, 5).Value = '  +0.04%  '
